$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.135.33'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.31%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.839.28'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.27%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.07'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6224'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.69%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.003'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.20%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07491'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.57%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2935'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.24'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.11%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07715'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.31%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.893.90'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.52%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.011'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6742'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.31%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.00'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.32%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009278'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -4.19%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.962'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.21%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.146.84'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.13%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.134.60'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.11%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '231.93'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.69'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.69%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.003'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.191'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.15%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.003'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '160.41'

$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.536'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.34%  '

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1390'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.76%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.88'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.50%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.166'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.90%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.134'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.89%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05567'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.25%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.35%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7522'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.74%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.845'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.08%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.142'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.06%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.663'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.10%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.769'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.37%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.223.94'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.91%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01786'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.34%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.519'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.74%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8990'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.99%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.002'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.14%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.024.34'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.57%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.97'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.11%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '66.08'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.58%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000121'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.89%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5097'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.23%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4085'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.11%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.097'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.68%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05843'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.18%  '
